$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New metric values (same for every model row, columns B..I)
$newVals = @(0.6731329884640765, 0.5777106424950826, -7.844374242436931, 0.2911304283046081, 0.3617455065250397, 0.4849758744239807, 0.4773751497268677, 0.4813994765281677)

# New ordering of the model labels for rows 2..26 (column A)
$newLabels = @(
  "model_10_7_0",
  "model_10_7_22",
  "model_10_7_21",
  "model_10_7_20",
  "model_10_7_19",
  "model_10_7_18",
  "model_10_7_17",
  "model_10_7_16",
  "model_10_7_15",
  "model_10_7_14",
  "model_10_7_13",
  "model_10_7_23",
  "model_10_7_12",
  "model_10_7_10",
  "model_10_7_9",
  "model_10_7_8",
  "model_10_7_7",
  "model_10_7_6",
  "model_10_7_5",
  "model_10_7_4",
  "model_10_7_3",
  "model_10_7_2",
  "model_10_7_1",
  "model_10_7_11",
  "model_10_7_24"
)

for ($i = 0; $i -lt $newLabels.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $newLabels[$i]
    for ($c = 0; $c -lt $newVals.Length; $c++) {
        $ws.Cells.Item($row, $c + 2).Value = $newVals[$c]
    }
}
